# ==========================================================================
# edit.ps1 - applies the "add faq feature, get all user, add notification
# isRead" commit to Architecture.docx via the Word COM object model.
# ==========================================================================

$d = $word.ActiveDocument
$ENDASH = [char]0x2013
$wdParagraph = 4

# --------------------------------------------------------------------------
# 1. "- User" section: add a new "- GET /user - Get all users" paragraph
#    right before the existing "- POST /user/register" paragraph.
#    (paragraph #20 == "- User" heading, in the original document)
# --------------------------------------------------------------------------
$userHeading = $d.Paragraphs.Item(20)   # "- User"
$pt = $d.Range($userHeading.Range.End, $userHeading.Range.End)
$pt.InsertAfter(" ")
$pt = $d.Range($pt.End, $pt.End)
$pt.InsertAfter(" ")
$pt = $d.Range($pt.End, $pt.End)
$pt.InsertAfter("- GET /user " + $ENDASH + " Get all users")
$pt = $d.Range($pt.End, $pt.End)
$pt.InsertParagraphAfter()

# --------------------------------------------------------------------------
# 2. Remove the stray "_GoBack" bookmark that used to sit on the
#    "- POST /deed/create" paragraph.
# --------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --------------------------------------------------------------------------
# 3. "- Log" section updates.
#    After step 1 we inserted exactly one new paragraph, so every paragraph
#    from #21 onward in the original document is now shifted down by 1.
#    "- POST /log/create ..." was #52 -> now #53.
#    "- GET /log/:id ..."     was #53 -> now #54.
# --------------------------------------------------------------------------
$logCreatePara = $d.Paragraphs.Item(53)

# 3a. Append the "(This is used for Share Copy link feature...)" note right
#     after "Create a message for the recipient." and before the line break.
$searchRng = $d.Range($logCreatePara.Range.Start, $logCreatePara.Range.End)
$searchRng.Find.Execute("Create a message for the recipient.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPt = $d.Range($searchRng.End, $searchRng.End)
$insPt.InsertAfter(" (This is used for Share Copy link feature. Users can send a request with deed_id)")

# 3b. Expand the Body Parameters list with "sender_name".
$logCreatePara = $d.Paragraphs.Item(53)
$searchRng2 = $d.Range($logCreatePara.Range.Start, $logCreatePara.Range.End)
$searchRng2.Find.Execute(": sender_id, recipient_id, message, ", $true, $false, $false, $false, $false, $true, 1, $false, ": sender_id, sender_name, recipient_id, message, ", 2) | Out-Null

# --------------------------------------------------------------------------
# 4. "- GET /log/:id" paragraph: append notification-setting note.
# --------------------------------------------------------------------------
$logGetPara = $d.Paragraphs.Item(54)
$endPt = $d.Range($logGetPara.Range.End - 1, $logGetPara.Range.End - 1)
$endPt.InsertAfter(". (This is used for notification setting)")

# --------------------------------------------------------------------------
# 5. Brand-new "- Faq" section appended after the Log section (paragraph 54
#    is still the last paragraph of the body, right before the sectPr).
# --------------------------------------------------------------------------
$logGetPara = $d.Paragraphs.Item(54)
$rng = $d.Range($logGetPara.Range.End, $logGetPara.Range.End)
$rng.InsertParagraphAfter()
$faqHeadingPara = $d.Paragraphs.Item(55)
$faqHeadingPara.Style = "Heading2"
$faqHeadingPara.Range.Text = "- Faq"

# -- "  - GET /faq  - GET ALL FAQs" ----------------------------------------
$faqHeadingPara = $d.Paragraphs.Item(55)
$rng = $d.Range($faqHeadingPara.Range.End, $faqHeadingPara.Range.End)
$rng.InsertParagraphAfter()
$getAllFaqPara = $d.Paragraphs.Item(56)
$getAllFaqPara.Range.Text = "  - GET /faq  - GET ALL FAQs"

# -- "  - POST /faq - Create a new FAQ" ------------------------------------
$getAllFaqPara = $d.Paragraphs.Item(56)
$rng = $d.Range($getAllFaqPara.Range.End, $getAllFaqPara.Range.End)
$rng.InsertParagraphAfter()
$createFaqPara = $d.Paragraphs.Item(57)
$createFaqPara.Range.Text = "  - POST /faq " + $ENDASH + " Create a new FAQ"

# -- "     Body Parameters: question, answer" (highlighted label) ---------
$createFaqPara = $d.Paragraphs.Item(57)
$rng = $d.Range($createFaqPara.Range.End, $createFaqPara.Range.End)
$rng.InsertParagraphAfter()
$bodyParams1Para = $d.Paragraphs.Item(58)
$insPt = $d.Range($bodyParams1Para.Range.Start, $bodyParams1Para.Range.Start)
$insPt.InsertAfter("     ")

# copy formatted "Body Parameters" (keeps the lightGray highlight) from an
# existing occurrence earlier in the document, then paste it in place.
$srcRng = $d.Content
$srcRng.Find.Execute("Body Parameters", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$srcRng.Copy()
$bodyParams1Para = $d.Paragraphs.Item(58)
$pastePt = $d.Range($bodyParams1Para.Range.End - 1, $bodyParams1Para.Range.End - 1)
$pastePt.Paste()

$bodyParams1Para = $d.Paragraphs.Item(58)
$pastePt2 = $d.Range($bodyParams1Para.Range.End - 1, $bodyParams1Para.Range.End - 1)
$pastePt2.InsertAfter(": question, answer")

# -- " - PATCH  /faq/:faq_id - Update an existing FAQ" ---------------------
$bodyParams1Para = $d.Paragraphs.Item(58)
$rng = $d.Range($bodyParams1Para.Range.End, $bodyParams1Para.Range.End)
$rng.InsertParagraphAfter()
$patchFaqPara = $d.Paragraphs.Item(59)
$patchFaqPara.Range.Text = " - PATCH  /faq/:faq_id " + $ENDASH + " Update an existing FAQ"

# -- "    Body Parameters: question, answer, status" (highlighted label) --
$patchFaqPara = $d.Paragraphs.Item(59)
$rng = $d.Range($patchFaqPara.Range.End, $patchFaqPara.Range.End)
$rng.InsertParagraphAfter()
$bodyParams2Para = $d.Paragraphs.Item(60)
$insPt = $d.Range($bodyParams2Para.Range.Start, $bodyParams2Para.Range.Start)
$insPt.InsertAfter("    ")

$srcRng2 = $d.Content
$srcRng2.Find.Execute("Body Parameters", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$srcRng2.Copy()
$bodyParams2Para = $d.Paragraphs.Item(60)
$pastePt3 = $d.Range($bodyParams2Para.Range.End - 1, $bodyParams2Para.Range.End - 1)
$pastePt3.Paste()

$bodyParams2Para = $d.Paragraphs.Item(60)
$pastePt4 = $d.Range($bodyParams2Para.Range.End - 1, $bodyParams2Para.Range.End - 1)
$pastePt4.InsertAfter(": question, answer, status")

# -- "- DELETE /faq/:faq_id - Delete an FAQ" + trailing _GoBack bookmark --
$bodyParams2Para = $d.Paragraphs.Item(60)
$rng = $d.Range($bodyParams2Para.Range.End, $bodyParams2Para.Range.End)
$rng.InsertParagraphAfter()
$deleteFaqPara = $d.Paragraphs.Item(61)
$deleteFaqPara.Range.Text = "- DELETE /faq/:faq_id " + $ENDASH + " Delete an FAQ"

$deleteFaqPara = $d.Paragraphs.Item(61)
$bmPt = $d.Range($deleteFaqPara.Range.End - 1, $deleteFaqPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmPt)

# -- trailing empty paragraph ----------------------------------------------
$deleteFaqPara = $d.Paragraphs.Item(61)
$rng = $d.Range($deleteFaqPara.Range.End, $deleteFaqPara.Range.End)
$rng.InsertParagraphAfter()

Write-Output "edit complete"
